$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.930.41'
$ws.Range("E2").Value = '  -2.23%  '

# Row 3
$ws.Range("D3").Value = '3.493.40'
$ws.Range("E3").Value = '  -1.23%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.33%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.47%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.97%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.633'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.52%  '

# Row 8
$ws.Range("E8").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.629'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.43%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.155'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.39%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.14'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.39%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000271'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.94%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.16'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.12%  '

# Row 14
$ws.Range("D14").Value = '4.049.20'
$ws.Range("E14").Value = '  -1.49%  '

# Row 15
$ws.Range("D15").Value = '3.490.97'
$ws.Range("E15").Value = '  -1.52%  '

# Row 16
$ws.Range("E16").Value = '  -0.36%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.28%  '

# Row 18
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.57%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '65.769.13'
$ws.Range("E19").Value = '  -2.59%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.998'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.68%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '412.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.89%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.55%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '85.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.15%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.19%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.75%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.99%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.83'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.46%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.98%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '30.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.64%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '617.49'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.86%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.35'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.34%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.37%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.109'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.24%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '59.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.95%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.151'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.09%  '

# Row 36
$ws.Range("E36").Value = '  +0.04%  '

# Row 37
$ws.Range("D37").Value = '0.0₃0795'
$ws.Range("E37").Value = '  -3.68%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.72%  '

# Row 39
$ws.Range("D39").Value = '3.309.47'
$ws.Range("E39").Value = '  +9.17%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.378'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.49%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.33%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.996'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.46%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.74%  '

# Row 44
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.03%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0413'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.04%  '

# Row 46
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.49%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.52%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.132'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.31%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '137.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.04%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.38'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.01%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.88%  '
